# Auto-generated Word COM-interop edit script
$d = $word.ActiveDocument

# Replace the (unique, literal) $oldText found inside paragraph $paraIndex with $newText.
# Always re-fetches the paragraph range first so Find stays scoped to that paragraph
# (a Find on a collapsed Range searches forward through the *whole* document, so we
# never reuse a stale collapsed Range across paragraph boundaries).
function Replace-InPara($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $ok = $p.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if (-not $ok) { throw "Replace-InPara: text not found in paragraph $($paraIndex): $oldText" }
    $p.Text = $newText
    return $newText
}

# Insert $newText as brand new run(s) right after the (unique, literal) $afterText
# inside paragraph $paraIndex. If $asBreak is $true, a manual line break (vertical
# tab / w:br) is inserted immediately before $newText, inside the new run, matching
# runs that pair a <w:br/> with following <w:t>.
function Insert-AfterAnchor($paraIndex, $afterText, $newText, [bool]$asBreak = $false) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $ok = $p.Find.Execute($afterText, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if (-not $ok) { throw "Insert-AfterAnchor: anchor not found in paragraph $($paraIndex): $afterText" }
    $p.Collapse(0) | Out-Null
    if ($asBreak) {
        $p.InsertAfter([char]11)
        $p.Collapse(0) | Out-Null
    }
    $p.InsertAfter($newText)
    return $newText
}

# --- Title / author / email (first three paragraphs) ---
Replace-InPara 1 'Orchestrating A Resilient Digital Domain' 'Venturing Through History''s Vast Tapestry' | Out-Null
Replace-InPara 2 'Liam Carson' 'Joshua Constantine' | Out-Null
Replace-InPara 3 'CarsonLi@wirednest' 'joshua' | Out-Null
Replace-InPara 3 'com' 'constantine@edusite' | Out-Null

# Append the remaining new runs onto the end of the e-mail address
Insert-AfterAnchor 3 'constantine@edusite' '.org' | Out-Null

# --- paragraph 5 ---
Replace-InPara 5 'In the contemporary era, the digital realm has emerged as a ubiquitous facet of human existence, profoundly influencing diverse aspects of our daily lives' 'History paints a vivid canvas of humanity''s triumphs, challenges and enduring legacies' | Out-Null

Replace-InPara 5 ' Its pervasiveness has sparked a critical need to ensure the resilience and security of this intricate cyber infrastructure' ' It''s a narrative that unfolds across diverse lands, centuries and cultures' | Out-Null

Replace-InPara 5 ' The protection of sensitive data, the integrity of digital systems, and the seamless functioning of critical services are paramount considerations in the face of ever-evolving cyber threats' ' We journey back in time to explore our roots, learn from past mistakes and gain perspective on the present' | Out-Null

Replace-InPara 5 ' As we navigate the complexities of the digital landscape, it is imperative that we prioritize the development and implementation of robust cybersecurity measures, fostering a resilient digital domain that can withstand the relentless barrage of malicious attacks' ' History not only satisfies our curiosity but also empowers us to navigate current events' | Out-Null

Insert-AfterAnchor 5 ' History not only satisfies our curiosity but also empowers us to navigate current events' '. The lives, ideas and events of yesterday shape our world today. By piecing together the historical narrative, we unravel the fabric of civilizations, appreciate the diversity of human experience and discover common threads that bind us across time' | Out-Null

Replace-InPara 5 'From financial transactions and confidential communications to intricate industrial control systems and intricate healthcare networks, the digital realm has become an indispensable component of modern society' 'Paragraph 2:' | Out-Null

Insert-AfterAnchor 5 'Paragraph 2:' 'History offers crucial lessons for the future' $true | Out-Null

Replace-InPara 5 ' However, this interconnectedness comes with inherent vulnerabilities, exposing individuals, organizations, and entire nations to a plethora of cyber risks' ' By delving into past conflicts, we learn the importance of peace and diplomacy' | Out-Null

Replace-InPara 5 ' Malicious actors, ranging from sophisticated state-sponsored groups to lone individuals with malicious intent, constantly exploit these vulnerabilities to perpetrate cyberattacks, leading to disruptions, data breaches, and even physical harm' ' Through studying economic downturns, we devise strategies to avert financial crises' | Out-Null

Replace-InPara 5 ' A resilient digital domain necessitates the adoption of a multi-pronged approach, encompassing technical safeguards such as encryption and multi-factor authentication, as well as robust cybersecurity policies and practices to mitigate these ever-present threats' ' The successes and failures of leaders past provide valuable insights for modern statecraft' | Out-Null

Insert-AfterAnchor 5 ' The successes and failures of leaders past provide valuable insights for modern statecraft' '. By examining societal and cultural transitions, we gain insight into the nuances of social evolution. History is a mirror that reflects our collective experiences, allowing us to recognize patterns, anticipate potential consequences and make informed decisions as individuals and societies' | Out-Null

Replace-InPara 5 'Furthermore, fostering a culture of cybersecurity awareness among all stakeholders is essential to cultivating a robust digital ecosystem' 'Paragraph 3:' | Out-Null

Insert-AfterAnchor 5 'Paragraph 3:' 'Unraveling history''s complexities requires diverse perspectives, critical thinking and an appreciation for nuance' $true | Out-Null

Replace-InPara 5 ' Equipping individuals with the knowledge to recognize and respond to cyber threats empowers them to become active participants in safeguarding the digital realm' ' History isn''t a linear, clear-cut narrative; it''s often messy, controversial and subject to interpretation' | Out-Null

Replace-InPara 5 ' Regular security audits and penetration testing can proactively identify vulnerabilities and weaknesses, enabling timely remediation measures to mitigate potential breaches' ' Understanding historical events demands an ability to weigh evidence, consider different viewpoints and engage in thoughtful discussion' | Out-Null

Replace-InPara 5 ' By proactively addressing these challenges, we can bolster the resilience of the digital domain and minimize the impact of malicious activities' ' Studying history helps cultivate these skills, fostering critical analysis and comprehensive understanding' | Out-Null

Insert-AfterAnchor 5 ' Studying history helps cultivate these skills, fostering critical analysis and comprehensive understanding' '. The study of history is an ongoing endeavor, with new discoveries and interpretations constantly emerging. Engaging with history requires curiosity, openness to new ideas and a willingness to challenge preconceptions' | Out-Null

# --- paragraph 7 ---
Replace-InPara 7 'The digital domain has become an integral part of our lives, necessitating the implementation of robust cybersecurity measures to ensure its resilience and protect against evolving cyber threats' 'History is a riveting tapestry of human experiences, offering profound lessons for the present and insights for the future' | Out-Null

Replace-InPara 7 ' A comprehensive approach encompassing technical safeguards, cybersecurity policies, and fostering awareness among stakeholders is essential to creating a cyber-resilient landscape' ' It calls us to explore diverse perspectives, think critically and appreciate the nuances of past events' | Out-Null

Replace-InPara 7 ' As we continue to leverage the benefits of digital technologies, prioritizing cybersecurity is paramount in safeguarding the integrity, confidentiality, and availability of critical digital systems and services' ' Through history, we unravel the enigmas of our origins, understand our collective triumphs and tribulations, and gain invaluable wisdom to navigate an ever-changing world' | Out-Null

Insert-AfterAnchor 7 ' Through history, we unravel the enigmas of our origins, understand our collective triumphs and tribulations, and gain invaluable wisdom to navigate an ever-changing world' '. Studying history is not just about memorizing names, dates and events; it''s about embarking on a journey through time, connecting with our ancestors and discovering the essence of what makes us human' | Out-Null

# --- Trailing empty paragraph appended at the very end of the document ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastPara.Collapse(0) | Out-Null
$lastPara.InsertParagraphAfter()

